$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 6 new rows before row 906 (shifts all existing rows 906:977 down to 912:983,
# growing the used range from A1:R977 to A1:R983). EntireRow.Insert() copies the
# formatting (incl. the date-number-format style on column D) from the row above,
# matching the existing layout.
$ws.Range("A906:A911").EntireRow.Insert()

# Populate the 6 newly inserted rows with this week's price report
# (Comercializadora del Agro de Limari - Pimiento, fecha 2022-08-10 / serial 44783).

# Row 906: Cuatro cascos verde - Primera
$ws.Cells.Item(906, 1).Value = 2
$ws.Cells.Item(906, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(906, 3).Value = "Coquimbo"
$ws.Cells.Item(906, 4).Value = 44783
$ws.Cells.Item(906, 5).Value = 4
$ws.Cells.Item(906, 6).Value = 100112002
$ws.Cells.Item(906, 7).Value = "Pimiento"
$ws.Cells.Item(906, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(906, 9).Value = "Primera"
$ws.Cells.Item(906, 10).Value = 1100
$ws.Cells.Item(906, 11).Value = 25000
$ws.Cells.Item(906, 12).Value = 26000
$ws.Cells.Item(906, 13).Value = 25500
$ws.Cells.Item(906, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(906, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(906, 16).Value = 1417
$ws.Cells.Item(906, 17).Value = 18
$ws.Cells.Item(906, 18).Value = "Hortaliza"

# Row 907: Cuatro cascos verde - Segunda
$ws.Cells.Item(907, 1).Value = 2
$ws.Cells.Item(907, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(907, 3).Value = "Coquimbo"
$ws.Cells.Item(907, 4).Value = 44783
$ws.Cells.Item(907, 5).Value = 4
$ws.Cells.Item(907, 6).Value = 100112002
$ws.Cells.Item(907, 7).Value = "Pimiento"
$ws.Cells.Item(907, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(907, 9).Value = "Segunda"
$ws.Cells.Item(907, 10).Value = 360
$ws.Cells.Item(907, 11).Value = 21000
$ws.Cells.Item(907, 12).Value = 22000
$ws.Cells.Item(907, 13).Value = 21500
$ws.Cells.Item(907, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(907, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(907, 16).Value = 1194
$ws.Cells.Item(907, 17).Value = 18
$ws.Cells.Item(907, 18).Value = "Hortaliza"

# Row 908: Cuatro cascos verde - Tercera
$ws.Cells.Item(908, 1).Value = 2
$ws.Cells.Item(908, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(908, 3).Value = "Coquimbo"
$ws.Cells.Item(908, 4).Value = 44783
$ws.Cells.Item(908, 5).Value = 4
$ws.Cells.Item(908, 6).Value = 100112002
$ws.Cells.Item(908, 7).Value = "Pimiento"
$ws.Cells.Item(908, 8).Value = "Cuatro cascos verde"
$ws.Cells.Item(908, 9).Value = "Tercera"
$ws.Cells.Item(908, 10).Value = 240
$ws.Cells.Item(908, 11).Value = 17000
$ws.Cells.Item(908, 12).Value = 18000
$ws.Cells.Item(908, 13).Value = 17500
$ws.Cells.Item(908, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(908, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(908, 16).Value = 972
$ws.Cells.Item(908, 17).Value = 18
$ws.Cells.Item(908, 18).Value = "Hortaliza"

# Row 909: Morrón rojo - Primera
$ws.Cells.Item(909, 1).Value = 2
$ws.Cells.Item(909, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(909, 3).Value = "Coquimbo"
$ws.Cells.Item(909, 4).Value = 44783
$ws.Cells.Item(909, 5).Value = 4
$ws.Cells.Item(909, 6).Value = 100112002
$ws.Cells.Item(909, 7).Value = "Pimiento"
$ws.Cells.Item(909, 8).Value = "Morrón rojo"
$ws.Cells.Item(909, 9).Value = "Primera"
$ws.Cells.Item(909, 10).Value = 600
$ws.Cells.Item(909, 11).Value = 27000
$ws.Cells.Item(909, 12).Value = 28000
$ws.Cells.Item(909, 13).Value = 27500
$ws.Cells.Item(909, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(909, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(909, 16).Value = 1528
$ws.Cells.Item(909, 17).Value = 18
$ws.Cells.Item(909, 18).Value = "Hortaliza"

# Row 910: Morrón rojo - Segunda
$ws.Cells.Item(910, 1).Value = 2
$ws.Cells.Item(910, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(910, 3).Value = "Coquimbo"
$ws.Cells.Item(910, 4).Value = 44783
$ws.Cells.Item(910, 5).Value = 4
$ws.Cells.Item(910, 6).Value = 100112002
$ws.Cells.Item(910, 7).Value = "Pimiento"
$ws.Cells.Item(910, 8).Value = "Morrón rojo"
$ws.Cells.Item(910, 9).Value = "Segunda"
$ws.Cells.Item(910, 10).Value = 400
$ws.Cells.Item(910, 11).Value = 23000
$ws.Cells.Item(910, 12).Value = 24000
$ws.Cells.Item(910, 13).Value = 23500
$ws.Cells.Item(910, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(910, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(910, 16).Value = 1306
$ws.Cells.Item(910, 17).Value = 18
$ws.Cells.Item(910, 18).Value = "Hortaliza"

# Row 911: Morrón rojo - Tercera
$ws.Cells.Item(911, 1).Value = 2
$ws.Cells.Item(911, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(911, 3).Value = "Coquimbo"
$ws.Cells.Item(911, 4).Value = 44783
$ws.Cells.Item(911, 5).Value = 4
$ws.Cells.Item(911, 6).Value = 100112002
$ws.Cells.Item(911, 7).Value = "Pimiento"
$ws.Cells.Item(911, 8).Value = "Morrón rojo"
$ws.Cells.Item(911, 9).Value = "Tercera"
$ws.Cells.Item(911, 10).Value = 340
$ws.Cells.Item(911, 11).Value = 19000
$ws.Cells.Item(911, 12).Value = 20000
$ws.Cells.Item(911, 13).Value = 19500
$ws.Cells.Item(911, 14).Value = "$/caja 18 kilos"
$ws.Cells.Item(911, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(911, 16).Value = 1083
$ws.Cells.Item(911, 17).Value = 18
$ws.Cells.Item(911, 18).Value = "Hortaliza"
